# Rótulos da primeira linha (cabeçalho) de cada planilha são alterados para
# incluir um prefixo textual ("Ano " ou "Intervalo "), de modo que o Power BI
# já reconheça automaticamente a primeira linha como cabeçalho da tabela.

$wb = $excel.ActiveWorkbook

# Planilhas cujo cabeçalho usa o prefixo "Ano "
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Planilha "Potencia Incremental - SIN(MW)" usa o prefixo "Intervalo "
$wsIncremental = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIncremental.Range("B1").Value = "Intervalo 2015"
$wsIncremental.Range("C1").Value = "Intervalo 2015-2030"
$wsIncremental.Range("D1").Value = "Intervalo 2031-2040"
$wsIncremental.Range("E1").Value = "Intervalo 2041-2050"

# Planilha "Custo Total (bilhões de R$)" só possui a coluna B no cabeçalho
$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value = "Ano 2015"
